# Turkey 1 Lig workbook update
# - Swap the full data (all columns except the leading sequence id in column A)
#   between several pairs of rows whose underlying fixtures were re-ordered.
# - Refresh the fixture data stored in rows 250-253 with newer odds/ids.
# - Remove the row that is no longer present (old row 254), whose content was
#   folded into row 250 (with refreshed odds) during the re-scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    $r1 = $ws.Range("B$row1`:AC$row1")
    $r2 = $ws.Range("B$row2`:AC$row2")
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value = $v2
    $r2.Value = $v1
}

# Pairs of rows whose full record (everything but the id in column A) is swapped
Swap-Rows 120 121
Swap-Rows 133 134
Swap-Rows 149 150
Swap-Rows 188 189
Swap-Rows 232 233

# Row 250: refreshed fixture (becomes what used to be the last row, with updated odds)
$ws.Range("B250").Value = 6963253
$ws.Range("E250").Value = 45386.60416666666
$ws.Range("F250").Value = "Erzurum BB"
$ws.Range("G250").Value = "Giresunspor"
$ws.Range("K250").Value = 1.142
$ws.Range("L250").Value = 6.5
$ws.Range("M250").Value = 15
$ws.Range("N250").Value = 1.166
$ws.Range("O250").Value = 6
$ws.Range("P250").Value = 17
$ws.Range("Q250").Value = -2
$ws.Range("T250").Value = 2.75
$ws.Range("U250").Value = 1.85
$ws.Range("V250").Value = 1.95

# Row 251: refreshed fixture
$ws.Range("B251").Value = 6963254
$ws.Range("E251").Value = 45388.3125
$ws.Range("F251").Value = "Bodrum BLD Spor"
$ws.Range("G251").Value = "Keciorengucu"
$ws.Range("K251").Value = 1.6
$ws.Range("L251").Value = 3.75
$ws.Range("M251").Value = 5.5
$ws.Range("N251").Value = 1.6
$ws.Range("O251").Value = 3.75
$ws.Range("P251").Value = 5.5
$ws.Range("Q251").Value = -0.75
$ws.Range("R251").Value = 1.775
$ws.Range("S251").Value = 2.025
$ws.Range("T251").Value = 2.25
$ws.Range("U251").Value = 1.8
$ws.Range("V251").Value = 2

# Row 252: refreshed fixture
$ws.Range("B252").Value = 6963050
$ws.Range("E252").Value = 45388.41666666666
$ws.Range("F252").Value = "Umraniyespor"
$ws.Range("G252").Value = "Sanliurfaspor"
$ws.Range("K252").Value = 2
$ws.Range("L252").Value = 3.4
$ws.Range("M252").Value = 3.5
$ws.Range("N252").Value = 2.375
$ws.Range("O252").Value = 3.25
$ws.Range("P252").Value = 2.8
$ws.Range("Q252").Value = -0.25
$ws.Range("R252").Value = 2.05
$ws.Range("S252").Value = 1.75
$ws.Range("T252").Value = 2.5
$ws.Range("U252").Value = 1.95
$ws.Range("V252").Value = 1.85

# Row 253: refreshed fixture
$ws.Range("B253").Value = 6963027
$ws.Range("E253").Value = 45388.60416666666
$ws.Range("F253").Value = "Sakaryaspor"
$ws.Range("G253").Value = "Boluspor"
$ws.Range("K253").Value = 1.85
$ws.Range("L253").Value = 3.25
$ws.Range("M253").Value = 4.333
$ws.Range("N253").Value = 1.85
$ws.Range("O253").Value = 3.25
$ws.Range("P253").Value = 4.333
$ws.Range("Q253").Value = -0.5
$ws.Range("R253").Value = 1.9
$ws.Range("S253").Value = 1.9
$ws.Range("T253").Value = 2.25
$ws.Range("U253").Value = 1.775
$ws.Range("V253").Value = 2.025

# Old row 254 is no longer part of the dataset; remove it entirely.
$ws.Rows.Item(254).Delete() | Out-Null
